# Adds two new bullet items right after the "a way to immediate remove
# from memory" bullet, and before the "Internal clone (...)" bullet:
#   - "optimize/review maintain_cache"                   (ilvl 0)
#   - "currently requiring all actions cleared before
#      maintain_cache, which is not optimal"             (ilvl 1)

$d = $word.ActiveDocument

# Locate the anchor paragraph via its known text.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*a way to immediate remove from memory*") {
        $anchor = $candidate
        break
    }
}

# Insert an empty paragraph right after the anchor, then fill it in with
# the exact OOXML for the "optimize/review maintain_cache" bullet (ilvl 0).
$anchor.Range.InsertParagraphAfter()
$firstNew = $anchor.Next()

$firstXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">optimize/review </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>maintain_cache</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$firstNew.Range.InsertXML($firstXml)

# Insert another paragraph after that one for the sub-bullet (ilvl 1)
# explaining why the optimization/review is needed.
$firstNew.Range.InsertParagraphAfter()
$secondNew = $firstNew.Next()

$secondXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">currently requiring all actions cleared before </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>maintain_cache</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>which is not optimal</w:t></w:r></w:p>
'@
$secondNew.Range.InsertXML($secondXml)
